# Generate Report for Handoff
# - Renames the tracked file from "a.md" -> the new GUID-based filename
#   (and drops the now-unused "b.md" row entirely).
# - Replaces the stale handoff/handback bookkeeping on the "zh-cn" and
#   "de-de" sheets with fresh handoff-only data (no handback yet), and
#   marks the old row as "Ignored" instead of "Include".
# - The ".localization-config" row shifts up into the vacated "b.md" row.

$wb = $excel.ActiveWorkbook

$newMdName  = "6403701d-b5be-4da2-bdab-d61aa05cb8b0.md"
$newMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/05f6ed0301a94a75ec73b781661b25f19c8d11aa/e2e/$newMdName"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/05f6ed0301a94a75ec73b781661b25f19c8d11aa/.localization-config"

$newZhXlfName = "6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.zh-cn.xlf"
$newZhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/919ef7093330fc6eceef7c1fd8224bab0d96b427/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlfName"

$newDeXlfName = "6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.de-de.xlf"
$newDeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4280317c23c4c9862d6eaa17028edfb4c25cf3b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlfName"

# ---------------------------------------------------------------------
# Sheet 1 "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newMdName

# Drop the "b.md" row (row 3); the ".localization-config" row (row 4)
# shifts up to become row 3.
$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, "", "", $newMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2 "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMdName
$ws2.Range("C2").Value = $newZhXlfName
$ws2.Range("D2").Value = "2016-03-11 05:55:54"
$ws2.Range("E2").Clear()
$ws2.Range("F2").Clear()
$ws2.Range("G2").Value = "0001-01-01 00:00:00"

# Drop the "b.md" row (row 3); the ".localization-config" row (row 4)
# shifts up to become row 3.
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, "", "", $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $newZhXlfUrl, "", "", $newZhXlfName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3 "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMdName
$ws3.Range("C2").Value = $newDeXlfName
$ws3.Range("D2").Value = "2016-03-11 05:57:39"
$ws3.Range("E2").Clear()
$ws3.Range("F2").Clear()
$ws3.Range("G2").Value = "0001-01-01 00:00:00"

# Drop the "b.md" row (row 3); the ".localization-config" row (row 4)
# shifts up to become row 3.
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, "", "", $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $newDeXlfUrl, "", "", $newDeXlfName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null
